# add rabbitmq for update data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the email shown in A2 (hyperlink target / style stay as-is)
$ws.Range("A2").Value = "acacacax21@yahoo.com"

# Row 3 previously held a second queued email + hyperlink; clear it down to a
# bare placeholder cell (keeps the Hyperlink style) and drop its hyperlink.
$ws.Range("A3:E3").ClearContents()

foreach ($hl in $ws.Hyperlinks) {
    if ($hl.Range.Address() -eq '$A$3') {
        $hl.Delete()
    }
}

# Two more blank placeholder rows, styled like the hyperlink rows, ready for
# the rabbitmq consumer to fill in.
$ws.Range("A4").Style = "Hyperlink"
$ws.Range("A5").Style = "Hyperlink"

# Widen column A so the longer addresses are readable (target stored width
# 32.21875; the engine rounds ColumnWidth to whole-pixel steps, so feed the
# input that lands closest to that stored value).
$ws.Columns.Item(1).ColumnWidth = 31.3

# Match the saved selection.
$null = $ws.Range("A2").Select()
